# Insert a new weekly data row at row 11 (pushing existing rows 11-48
# down to 12-49), then populate the new row 11 with its values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 11..48 down to 12..49 (equivalent to Excel's right-click -> Insert on row 11)
$ws.Rows.Item(11).Insert()

# Populate the newly inserted row 11
$ws.Cells.Item(11, 1).Value = 11
$ws.Cells.Item(11, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(11, 3).Value = "Bíobío"
$ws.Cells.Item(11, 4).Value = 45035
$ws.Cells.Item(11, 5).Value = 8
$ws.Cells.Item(11, 6).Value = 100112026
$ws.Cells.Item(11, 7).Value = "Haba"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 100
$ws.Cells.Item(11, 11).Value = 18000
$ws.Cells.Item(11, 12).Value = 20000
$ws.Cells.Item(11, 13).Value = 19000
$ws.Cells.Item(11, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(11, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(11, 16).Value = 760
$ws.Cells.Item(11, 17).Value = 25
$ws.Cells.Item(11, 18).Value = "Hortaliza"

# Make sure the style for the date column D matches the other date cells
$ws.Cells.Item(11, 4).NumberFormat = $ws.Cells.Item(12, 4).NumberFormat
